$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseDataSets")

$ws.Range("B2").Interior.Color = 5296274
$ws.Range("B3").Interior.Color = 5296274

$ws.Range("G41").Interior.Color = 16777215
Write-Host "done"
